# Refresh "想去人数" (F column, want-to-go count) values on the 展览 and
# 全部类型 sheets to match the newly generated output.

$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 4665
$ws1.Cells.Item(3, 6).Value = 2538
$ws1.Cells.Item(8, 6).Value = 235
$ws1.Cells.Item(10, 6).Value = 195
$ws1.Cells.Item(11, 6).Value = 179
$ws1.Cells.Item(12, 6).Value = 1747
$ws1.Cells.Item(14, 6).Value = 3903
$ws1.Cells.Item(15, 6).Value = 38
$ws1.Cells.Item(16, 6).Value = 263

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 4665
$ws4.Cells.Item(3, 6).Value = 2538
$ws4.Cells.Item(10, 6).Value = 235
$ws4.Cells.Item(12, 6).Value = 195
$ws4.Cells.Item(13, 6).Value = 179
$ws4.Cells.Item(16, 6).Value = 1747
$ws4.Cells.Item(18, 6).Value = 3903
$ws4.Cells.Item(19, 6).Value = 38
$ws4.Cells.Item(20, 6).Value = 263
